# Scheduled Universalis price refresh for the Goblin_Profits workbook.
# Re-pulls currentAveragePrice(NQ/HQ) for each leve whose market price
# moved, then recomputes the dependent Leve*Price*/Leve*Profit* columns.
# (Columns: H currentAveragePrice, I currentAveragePriceNQ,
#  J currentAveragePriceHQ, K LevePriceNQ, L LevePriceHQ,
#  M LeveProfitNQ, N LeveProfitHQ.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 71.333336
$ws.Range("I5").Value = 44.857143
$ws.Range("J5").Value = 164
$ws.Range("K5").Value = 44.857143
$ws.Range("L5").Value = 164
$ws.Range("M5").Value = 70.14285699999999
$ws.Range("N5").Value = -394
$ws.Range("H17").Value = 1055.878
$ws.Range("J17").Value = 1077.55
$ws.Range("L17").Value = 3232.65
$ws.Range("N17").Value = -3568.65
$ws.Range("H18").Value = 2425
$ws.Range("I18").Value = 566.6667
$ws.Range("K18").Value = 566.6667
$ws.Range("M18").Value = -282.6667
$ws.Range("H28").Value = 5006.696
$ws.Range("I28").Value = 6484.8237
$ws.Range("K28").Value = 6484.8237
$ws.Range("M28").Value = -5999.8237
$ws.Range("H40").Value = 2380.524
$ws.Range("J40").Value = 2999.75
$ws.Range("L40").Value = 2999.75
$ws.Range("N40").Value = -3349.75
$ws.Range("H58").Value = 25001158
$ws.Range("I58").Value = 50000320
$ws.Range("J58").Value = 1997.2
$ws.Range("K58").Value = 150000960
$ws.Range("L58").Value = 5991.6
$ws.Range("M58").Value = -150000810
$ws.Range("N58").Value = -6291.6
$ws.Range("H62").Value = 57988
$ws.Range("I62").Value = 103663.5
$ws.Range("J62").Value = 12312.5
$ws.Range("K62").Value = 103663.5
$ws.Range("L62").Value = 12312.5
$ws.Range("M62").Value = -103039.5
$ws.Range("N62").Value = -13560.5
$ws.Range("H65").Value = 57988
$ws.Range("I65").Value = 103663.5
$ws.Range("J65").Value = 12312.5
$ws.Range("K65").Value = 518317.5
$ws.Range("L65").Value = 61562.5
$ws.Range("M65").Value = -515197.5
$ws.Range("N65").Value = -67802.5
$ws.Range("H74").Value = 4107.8184
$ws.Range("I74").Value = 4107.8184
$ws.Range("K74").Value = 4107.8184
$ws.Range("M74").Value = -3171.8184
$ws.Range("H77").Value = 4107.8184
$ws.Range("I77").Value = 4107.8184
$ws.Range("K77").Value = 20539.092
$ws.Range("M77").Value = -15859.092
$ws.Range("H80").Value = 527.2105
$ws.Range("I80").Value = 128.5
$ws.Range("J80").Value = 817.1818
$ws.Range("K80").Value = 385.5
$ws.Range("L80").Value = 2451.5454
$ws.Range("M80").Value = 612.5
$ws.Range("N80").Value = -4447.5454
$ws.Range("H83").Value = 527.2105
$ws.Range("I83").Value = 128.5
$ws.Range("J83").Value = 817.1818
$ws.Range("K83").Value = 1156.5
$ws.Range("L83").Value = 7354.6362
$ws.Range("M83").Value = 3835.5
$ws.Range("N83").Value = -17338.6362
$ws.Range("H98").Value = 12928.267
$ws.Range("I98").Value = 17593.7
$ws.Range("K98").Value = 17593.7
$ws.Range("M98").Value = -16095.7
$ws.Range("H113").Value = 3564.6072
$ws.Range("I113").Value = 3339.3044
$ws.Range("J113").Value = 4601
$ws.Range("K113").Value = 3339.3044
$ws.Range("L113").Value = 4601
$ws.Range("M113").Value = -85.30439999999999
$ws.Range("N113").Value = -11109
$ws.Range("H121").Value = 2501
$ws.Range("J121").Value = 2501
$ws.Range("L121").Value = 7503
$ws.Range("N121").Value = -10997
$ws.Range("H122").Value = 12928.267
$ws.Range("I122").Value = 17593.7
$ws.Range("K122").Value = 52781.10000000001
$ws.Range("M122").Value = -50331.10000000001
$ws.Range("H129").Value = 2194.111
$ws.Range("I129").Value = 949.4
$ws.Range("J129").Value = 3750
$ws.Range("K129").Value = 2848.2
$ws.Range("L129").Value = 11250
$ws.Range("M129").Value = 2151.8
$ws.Range("N129").Value = -21250
$ws.Range("H131").Value = 7344.08
$ws.Range("I131").Value = 5906.933
$ws.Range("J131").Value = 9499.799999999999
$ws.Range("K131").Value = 17720.799
$ws.Range("L131").Value = 28499.4
$ws.Range("M131").Value = -12680.799
$ws.Range("N131").Value = -38579.39999999999
$ws.Range("H132").Value = 1595.0571
$ws.Range("I132").Value = 1346.1562
$ws.Range("K132").Value = 4038.4686
$ws.Range("M132").Value = -1508.4686
$ws.Range("H138").Value = 2282.4807
$ws.Range("I138").Value = 1222.125
$ws.Range("J138").Value = 2753.75
$ws.Range("K138").Value = 3666.375
$ws.Range("L138").Value = 8261.25
$ws.Range("M138").Value = 1473.625
$ws.Range("N138").Value = -18541.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2605.9375
$ws.Range("I2").Value = 1169.6
$ws.Range("J2").Value = 4999.8335
$ws.Range("K2").Value = 1169.6
$ws.Range("L2").Value = 4999.8335
$ws.Range("M2").Value = -1056.6
$ws.Range("N2").Value = -5225.8335
$ws.Range("H32").Value = 3462.4119
$ws.Range("I32").Value = 3462.4119
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3462.4119
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3175.4119
$ws.Range("N32").ClearContents()
$ws.Range("H38").Value = 4748.375
$ws.Range("I38").Value = 4580.857
$ws.Range("K38").Value = 4580.857
$ws.Range("M38").Value = -4113.857
$ws.Range("H45").Value = 1869.2858
$ws.Range("I45").Value = 1597.5
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 1597.5
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -1220.5
$ws.Range("N45").Value = -4254
$ws.Range("H61").Value = 5354.316
$ws.Range("I61").Value = 5916.3335
$ws.Range("K61").Value = 5916.3335
$ws.Range("M61").Value = -5704.3335
$ws.Range("H104").Value = 3204
$ws.Range("I104").Value = 3204
$ws.Range("K104").Value = 3204
$ws.Range("M104").Value = 290
$ws.Range("H116").Value = 2605.9375
$ws.Range("I116").Value = 1169.6
$ws.Range("J116").Value = 4999.8335
$ws.Range("K116").Value = 1169.6
$ws.Range("L116").Value = 4999.8335
$ws.Range("M116").Value = 1124.4
$ws.Range("N116").Value = -9587.833500000001
$ws.Range("H122").Value = 1391.3903
$ws.Range("I122").Value = 1229.6389
$ws.Range("J122").Value = 2556
$ws.Range("K122").Value = 3688.9167
$ws.Range("L122").Value = 7668
$ws.Range("M122").Value = -1238.9167
$ws.Range("N122").Value = -12568
$ws.Range("H132").Value = 3406
$ws.Range("I132").Value = 2537.4546
$ws.Range("K132").Value = 7612.3638
$ws.Range("M132").Value = -5082.3638
$ws.Range("H136").Value = 5354.316
$ws.Range("I136").Value = 5916.3335
$ws.Range("K136").Value = 17749.0005
$ws.Range("M136").Value = -15199.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2605.9375
$ws.Range("I3").Value = 1169.6
$ws.Range("J3").Value = 4999.8335
$ws.Range("K3").Value = 1169.6
$ws.Range("L3").Value = 4999.8335
$ws.Range("M3").Value = -1055.6
$ws.Range("N3").Value = -5227.8335
$ws.Range("H20").Value = 1139
$ws.Range("I20").Value = 1400
$ws.Range("J20").Value = 1008.5
$ws.Range("K20").Value = 1400
$ws.Range("L20").Value = 1008.5
$ws.Range("M20").Value = -1153
$ws.Range("N20").Value = -1502.5
$ws.Range("H22").Value = 1432.6428
$ws.Range("I22").Value = 1265.4
$ws.Range("K22").Value = 1265.4
$ws.Range("M22").Value = -1092.4
$ws.Range("H61").Value = 75000
$ws.Range("J61").Value = 75000
$ws.Range("L61").Value = 75000
$ws.Range("N61").Value = -75626
$ws.Range("H80").Value = 378.76
$ws.Range("I80").Value = 382.16666
$ws.Range("J80").Value = 377.6842
$ws.Range("K80").Value = 382.16666
$ws.Range("L80").Value = 377.6842
$ws.Range("M80").Value = 615.83334
$ws.Range("N80").Value = -2373.6842
$ws.Range("H83").Value = 378.76
$ws.Range("I83").Value = 382.16666
$ws.Range("J83").Value = 377.6842
$ws.Range("K83").Value = 1910.8333
$ws.Range("L83").Value = 1888.421
$ws.Range("M83").Value = 3081.1667
$ws.Range("N83").Value = -11872.421
$ws.Range("H86").Value = 3137.1765
$ws.Range("I86").Value = 2757.9167
$ws.Range("J86").Value = 4047.4
$ws.Range("K86").Value = 2757.9167
$ws.Range("L86").Value = 4047.4
$ws.Range("M86").Value = -1634.9167
$ws.Range("N86").Value = -6293.4
$ws.Range("H89").Value = 3137.1765
$ws.Range("I89").Value = 2757.9167
$ws.Range("J89").Value = 4047.4
$ws.Range("K89").Value = 13789.5835
$ws.Range("L89").Value = 20237
$ws.Range("M89").Value = -8173.583500000001
$ws.Range("N89").Value = -31469
$ws.Range("H94").Value = 2460.0715
$ws.Range("I94").Value = 1999.8
$ws.Range("K94").Value = 1999.8
$ws.Range("M94").Value = -1548.8
$ws.Range("H105").Value = 2909.9167
$ws.Range("J105").Value = 2622.2
$ws.Range("L105").Value = 2622.2
$ws.Range("N105").Value = -6116.2
$ws.Range("H123").Value = 105000
$ws.Range("J123").Value = 105000
$ws.Range("L123").Value = 105000
$ws.Range("N123").Value = -114800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3444.8333
$ws.Range("I16").Value = 3406
$ws.Range("K16").Value = 3406
$ws.Range("M16").Value = -3119
$ws.Range("H22").Value = 1130.1818
$ws.Range("I22").Value = 815.7619
$ws.Range("J22").Value = 1680.4166
$ws.Range("K22").Value = 815.7619
$ws.Range("L22").Value = 1680.4166
$ws.Range("M22").Value = -465.7619
$ws.Range("N22").Value = -2380.4166
$ws.Range("H28").Value = 62000
$ws.Range("J28").Value = 62000
$ws.Range("L28").Value = 62000
$ws.Range("N28").Value = -62490
$ws.Range("H31").Value = 6533.0835
$ws.Range("I31").Value = 2050
$ws.Range("K31").Value = 2050
$ws.Range("M31").Value = -1755
$ws.Range("H34").Value = 6533.0835
$ws.Range("I34").Value = 2050
$ws.Range("K34").Value = 2050
$ws.Range("M34").Value = -1848
$ws.Range("H58").Value = 1379.6471
$ws.Range("I58").Value = 1343.0834
$ws.Range("J58").Value = 1467.4
$ws.Range("K58").Value = 1343.0834
$ws.Range("L58").Value = 1467.4
$ws.Range("M58").Value = -1140.0834
$ws.Range("N58").Value = -1873.4
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H94").Value = 3218.125
$ws.Range("I94").Value = 2660
$ws.Range("J94").Value = 4148.3335
$ws.Range("K94").Value = 2660
$ws.Range("L94").Value = 4148.3335
$ws.Range("M94").Value = -2209
$ws.Range("N94").Value = -5050.3335
$ws.Range("H95").Value = 6700
$ws.Range("J95").Value = 6700
$ws.Range("L95").Value = 6700
$ws.Range("N95").Value = -12192
$ws.Range("H99").Value = 3054
$ws.Range("I99").Value = 3054
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3054
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1556
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 3444.8333
$ws.Range("I113").Value = 3406
$ws.Range("K113").Value = 3406
$ws.Range("M113").Value = -1236
$ws.Range("H122").Value = 1729.1666
$ws.Range("I122").Value = 1626.6666
$ws.Range("J122").Value = 1831.6666
$ws.Range("K122").Value = 4879.9998
$ws.Range("L122").Value = 5494.9998
$ws.Range("M122").Value = -2429.9998
$ws.Range("N122").Value = -10394.9998
$ws.Range("H126").Value = 3054
$ws.Range("I126").Value = 3054
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9162
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6692
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3115.75
$ws.Range("I132").Value = 3221
$ws.Range("J132").Value = 2884.2
$ws.Range("K132").Value = 9663
$ws.Range("L132").Value = 8652.599999999999
$ws.Range("M132").Value = -7133
$ws.Range("N132").Value = -13712.6
$ws.Range("H134").Value = 3378.2144
$ws.Range("I134").Value = 4817
$ws.Range("K134").Value = 14451
$ws.Range("M134").Value = -11916
$ws.Range("H136").Value = 1379.6471
$ws.Range("I136").Value = 1343.0834
$ws.Range("J136").Value = 1467.4
$ws.Range("K136").Value = 4029.2502
$ws.Range("L136").Value = 4402.200000000001
$ws.Range("M136").Value = -1479.2502
$ws.Range("N136").Value = -9502.200000000001
$ws.Range("H141").Value = 148407.25
$ws.Range("J141").Value = 148407.25
$ws.Range("L141").Value = 148407.25
$ws.Range("N141").Value = -158767.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 96.63158
$ws.Range("J12").Value = 107.882355
$ws.Range("L12").Value = 323.647065
$ws.Range("N12").Value = -669.647065
$ws.Range("H32").Value = 1000000
$ws.Range("J32").Value = 1000000
$ws.Range("L32").Value = 3000000
$ws.Range("N32").Value = -3000566
$ws.Range("H33").Value = 162.75
$ws.Range("I33").Value = 50
$ws.Range("J33").Value = 200.33333
$ws.Range("K33").Value = 300
$ws.Range("L33").Value = 1201.99998
$ws.Range("M33").Value = -17
$ws.Range("N33").Value = -1767.99998
$ws.Range("H68").Value = 4610.1113
$ws.Range("J68").Value = 5061.5
$ws.Range("L68").Value = 15184.5
$ws.Range("N68").Value = -16806.5
$ws.Range("H71").Value = 4610.1113
$ws.Range("J71").Value = 5061.5
$ws.Range("L71").Value = 45553.5
$ws.Range("N71").Value = -53665.5
$ws.Range("H107").Value = 2780.6365
$ws.Range("J107").Value = 1198.25
$ws.Range("L107").Value = 3594.75
$ws.Range("N107").Value = -7434.75
$ws.Range("H113").Value = 1257.9286
$ws.Range("I113").Value = 725
$ws.Range("J113").Value = 1346.75
$ws.Range("K113").Value = 2175
$ws.Range("L113").Value = 4040.25
$ws.Range("M113").Value = -5
$ws.Range("N113").Value = -8380.25
$ws.Range("H120").Value = 113834.8
$ws.Range("I120").Value = 343844
$ws.Range("K120").Value = 1031532
$ws.Range("M120").Value = -1026694
$ws.Range("H121").Value = 100680.3
$ws.Range("I121").Value = 143400.42
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 430201.26
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = -428891.26
$ws.Range("N121").Value = -5620
$ws.Range("H122").Value = 2444.0557
$ws.Range("J122").Value = 2874.1538
$ws.Range("L122").Value = 25867.3842
$ws.Range("N122").Value = -30767.3842
$ws.Range("H132").Value = 2228.2144
$ws.Range("J132").Value = 2999.8572
$ws.Range("L132").Value = 26998.7148
$ws.Range("N132").Value = -32058.7148
$ws.Range("H140").Value = 53864
$ws.Range("I140").Value = 59730.35
$ws.Range("K140").Value = 179191.05
$ws.Range("M140").Value = -174011.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 938.65216
$ws.Range("I2").Value = 1709
$ws.Range("J2").Value = 98.27273
$ws.Range("K2").Value = 1709
$ws.Range("L2").Value = 98.27273
$ws.Range("M2").Value = -1596
$ws.Range("N2").Value = -324.27273
$ws.Range("H49").Value = 23801.143
$ws.Range("J49").Value = 23801.143
$ws.Range("L49").Value = 23801.143
$ws.Range("N49").Value = -24169.143
$ws.Range("H70").Value = 15372.25
$ws.Range("I70").Value = 10600.6
$ws.Range("K70").Value = 10600.6
$ws.Range("M70").Value = -10330.6
$ws.Range("H73").Value = 15372.25
$ws.Range("I73").Value = 10600.6
$ws.Range("K73").Value = 10600.6
$ws.Range("M73").Value = -9664.6
$ws.Range("H97").Value = 6387.2354
$ws.Range("I97").Value = 623.8182
$ws.Range("K97").Value = 623.8182
$ws.Range("M97").Value = -127.8182
$ws.Range("H104").Value = 198999
$ws.Range("J104").Value = 198999
$ws.Range("L104").Value = 198999
$ws.Range("N104").Value = -205987
$ws.Range("H105").Value = 31500
$ws.Range("J105").Value = 31500
$ws.Range("L105").Value = 31500
$ws.Range("N105").Value = -38488
$ws.Range("H106").Value = 49279.5
$ws.Range("J106").Value = 49279.5
$ws.Range("L106").Value = 49279.5
$ws.Range("N106").Value = -51803.5
$ws.Range("H113").Value = 5351.36
$ws.Range("I113").Value = 1543.9231
$ws.Range("K113").Value = 1543.9231
$ws.Range("M113").Value = 626.0769
$ws.Range("H122").Value = 5314.1875
$ws.Range("I122").Value = 5036.3794
$ws.Range("J122").Value = 7999.6665
$ws.Range("K122").Value = 15109.1382
$ws.Range("L122").Value = 23998.9995
$ws.Range("M122").Value = -12659.1382
$ws.Range("N122").Value = -28898.9995
$ws.Range("H132").Value = 2326.6206
$ws.Range("I132").Value = 2513.0476
$ws.Range("K132").Value = 7539.1428
$ws.Range("M132").Value = -5009.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3976.182
$ws.Range("I7").Value = 3988.8
$ws.Range("J7").Value = 3850
$ws.Range("K7").Value = 3988.8
$ws.Range("L7").Value = 3850
$ws.Range("M7").Value = -3876.8
$ws.Range("N7").Value = -4074
$ws.Range("H16").Value = 621.875
$ws.Range("I16").Value = 735
$ws.Range("K16").Value = 735
$ws.Range("M16").Value = -565
$ws.Range("H22").Value = 2865
$ws.Range("I22").Value = 2326.2778
$ws.Range("J22").Value = 3375.3684
$ws.Range("K22").Value = 2326.2778
$ws.Range("L22").Value = 3375.3684
$ws.Range("M22").Value = -2031.2778
$ws.Range("N22").Value = -3965.3684
$ws.Range("H27").Value = 2865
$ws.Range("I27").Value = 2326.2778
$ws.Range("J27").Value = 3375.3684
$ws.Range("K27").Value = 2326.2778
$ws.Range("L27").Value = 3375.3684
$ws.Range("M27").Value = -2219.2778
$ws.Range("N27").Value = -3589.3684
$ws.Range("H35").Value = 3500
$ws.Range("I35").Value = 3500
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3500
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -3164
$ws.Range("N35").ClearContents()
$ws.Range("H46").Value = 2297.111
$ws.Range("I46").Value = 893.3333
$ws.Range("K46").Value = 893.3333
$ws.Range("M46").Value = -705.3333
$ws.Range("H122").Value = 4163.25
$ws.Range("I122").Value = 3954.8667
$ws.Range("K122").Value = 11864.6001
$ws.Range("M122").Value = -9414.6001
$ws.Range("H126").Value = 3976.182
$ws.Range("I126").Value = 3988.8
$ws.Range("J126").Value = 3850
$ws.Range("K126").Value = 11966.4
$ws.Range("L126").Value = 11550
$ws.Range("M126").Value = -9496.400000000001
$ws.Range("N126").Value = -16490
$ws.Range("H136").Value = 5221.4443
$ws.Range("I136").Value = 5284.857
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 15854.571
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -13304.571
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2600
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 2600
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H103").Value = 34950
$ws.Range("J103").Value = 34950
$ws.Range("L103").Value = 34950
$ws.Range("N103").Value = -37294
$ws.Range("H107").Value = 4312.5835
$ws.Range("I107").Value = 2284.9443
$ws.Range("J107").Value = 10395.5
$ws.Range("K107").Value = 6854.8329
$ws.Range("L107").Value = 31186.5
$ws.Range("M107").Value = -4934.8329
$ws.Range("N107").Value = -35026.5
$ws.Range("H113").Value = 1819.7916
$ws.Range("I113").Value = 1727.0588
$ws.Range("J113").Value = 2045
$ws.Range("K113").Value = 5181.1764
$ws.Range("L113").Value = 6135
$ws.Range("M113").Value = -3011.1764
$ws.Range("N113").Value = -10475
$ws.Range("H122").Value = 6059.3125
$ws.Range("I122").Value = 1278.5714
$ws.Range("J122").Value = 9777.666999999999
$ws.Range("K122").Value = 3835.7142
$ws.Range("L122").Value = 29333.001
$ws.Range("M122").Value = -1385.7142
$ws.Range("N122").Value = -34233.001
$ws.Range("H126").Value = 1414.909
$ws.Range("I126").Value = 1096.125
$ws.Range("K126").Value = 3288.375
$ws.Range("M126").Value = -818.375
$ws.Range("H132").Value = 4456.9585
$ws.Range("I132").Value = 4433.1304
$ws.Range("K132").Value = 13299.3912
$ws.Range("M132").Value = -10769.3912
$ws.Range("H136").Value = 1651.9048
$ws.Range("J136").Value = 2909.5
$ws.Range("L136").Value = 8728.5
$ws.Range("N136").Value = -13828.5
